$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# CU - Registrar asistencia (row 6): status changes from "Por iniciar" to "En proceso"
$ws.Range("F6").Value = "En proceso"

# CU - CRU renta de espacio (row 8): status changes from "En proceso" to "Hecho"
$ws.Range("F8").Value = "Hecho"

# Día 3 consumption (column N): 1 hour consumed for both tasks
$ws.Range("N6").Value = 1
$ws.Range("N8").Value = 1

# Update the active selection to reflect where the user clicked last (O10)
$ws.Range("O10").Select()
